$wb = $excel.ActiveWorkbook

# --- Sheet: Summary ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B2").Value = 0.5
$summary.Range("C2").Value = 0.5
$summary.Range("D2").Value = 1
$summary.Range("E2").Value = 0.6666666666666666
$summary.Range("F2").Value = 0.8333333333333334
$summary.Range("G2").Value = 0.9629629629629629
$summary.Range("H2").Value = 0.7901113776318927
$summary.Range("I2").Value = 534
$summary.Range("J2").Value = 534
$summary.Range("K2").Value = 0
$summary.Range("L2").Value = 0

# --- Sheet: Classification Report ---
$clf = $wb.Worksheets.Item("Classification Report")

# row 2 (label "0")
$clf.Range("B2").Value = 0
$clf.Range("C2").Value = 0
$clf.Range("D2").Value = 0

# row 3 (label "1")
$clf.Range("B3").Value = 0.5
$clf.Range("C3").Value = 1
$clf.Range("D3").Value = 0.6666666666666666

# row 4 (accuracy)
$clf.Range("B4").Value = 0.5
$clf.Range("C4").Value = 0.5
$clf.Range("D4").Value = 0.5
$clf.Range("E4").Value = 0.5

# row 5 (macro avg)
$clf.Range("B5").Value = 0.25
$clf.Range("C5").Value = 0.5
$clf.Range("D5").Value = 0.3333333333333333

# row 6 (weighted avg)
$clf.Range("B6").Value = 0.25
$clf.Range("C6").Value = 0.5
$clf.Range("D6").Value = 0.3333333333333333

# --- Sheet: Confusion Matrix ---
$cm = $wb.Worksheets.Item("Confusion Matrix")
$cm.Range("B2").Value = 0
$cm.Range("C2").Value = 534
$cm.Range("B3").Value = 0
$cm.Range("C3").Value = 534
